# Challenge 4: Redone charts
# - Drop the rows for states/regions that are no longer present in the
#   final dataset (BW, CA, ELN, FIF, KHL, MA, MN, NW), which shifts the
#   remaining 11 rows (AZ, EDH, IL, MLN, NC, NV, ON, PA, QC, SC, WI) up to
#   occupy rows 2-12.
# - Re-enter the %reviews-quality formula as a single range fill so it
#   collapses into one shared formula (D2:D12).
# - Point both charts' series at the new, smaller ranges.
# - Give the stacked-bar chart a "#reviews" title.
# - Mark the line chart's value axis as a 0.00% percentage format.
# - Resize/reposition both chart frames to their new extents.
# - Restore the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the rows that disappear from the dataset. Delete from the
#    bottom up so earlier row numbers stay valid while iterating.
# ---------------------------------------------------------------------
$ws.Rows.Item(15).Delete()   # NW
$ws.Rows.Item(12).Delete()   # MN
$ws.Rows.Item(10).Delete()   # MA
$ws.Rows.Item(9).Delete()    # KHL
$ws.Rows.Item(7).Delete()    # FIF
$ws.Rows.Item(6).Delete()    # ELN
$ws.Rows.Item(4).Delete()    # CA
$ws.Rows.Item(3).Delete()    # BW

# ---------------------------------------------------------------------
# 2. Refill the ratio column as a single range formula so it is stored
#    as one shared formula spanning D2:D12.
# ---------------------------------------------------------------------
$ws.Range("D2:D12").Formula = "=B2/C2"

# ---------------------------------------------------------------------
# 3. Point the chart series at the new 11-row ranges.
# ---------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$chart1 = $co1.Chart
$bar1 = $chart1.SeriesCollection(1)
$bar1.Formula = '=SERIES(Sheet1!$B$1,Sheet1!$A$2:$A$12,Sheet1!$B$2:$B$12,1)'
$bar2 = $chart1.SeriesCollection(2)
$bar2.Formula = '=SERIES(Sheet1!$C$1,Sheet1!$A$2:$A$12,Sheet1!$C$2:$C$12,2)'

$co2 = $ws.ChartObjects().Item(2)
$chart2 = $co2.Chart
$line1 = $chart2.SeriesCollection(1)
$line1.Formula = '=SERIES(Sheet1!$D$1,Sheet1!$A$2:$A$12,Sheet1!$D$2:$D$12,1)'

# ---------------------------------------------------------------------
# 4. Add the "#reviews" title to the stacked bar chart.
# ---------------------------------------------------------------------
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "#reviews"

# ---------------------------------------------------------------------
# 5. Mark the line chart's value axis as a percentage format.
# ---------------------------------------------------------------------
$valAx2 = $chart2.Axes(2)
$valAx2.NumberFormatLinked = $false
$valAx2.NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 6. Resize/reposition the chart frames. Compute the deltas from the
#    actual column widths / row heights so the new edges line up with
#    the same cell boundaries the reference workbook uses.
# ---------------------------------------------------------------------
$EMU_PER_PT = 12700

function ColsWidthPt($c1, $c2) {
    $sum = 0
    for ($c = $c1; $c -le $c2; $c++) {
        $sum = $sum + $ws.Columns.Item($c).Width
    }
    return $sum
}
function RowsHeightPt($r1, $r2) {
    $sum = 0
    for ($r = $r1; $r -le $r2; $r++) {
        $sum = $sum + $ws.Rows.Item($r).Height
    }
    return $sum
}

# Chart 1 ("to" corner moves from col26/row42 to col30/row41, 1-based):
# grows right across columns 26-29, shrinks up across row 41.
$chart1ColGain = ColsWidthPt 26 29
$chart1OffXDelta = 43252 - 209550
$chart1OffXDeltaPt = $chart1OffXDelta / $EMU_PER_PT
$chart1WidthDelta = $chart1ColGain + $chart1OffXDeltaPt

$chart1RowLoss = RowsHeightPt 41 41
$chart1OffYDelta = 63500 - 177800
$chart1OffYDeltaPt = $chart1OffYDelta / $EMU_PER_PT
$chart1HeightDelta = (0 - $chart1RowLoss) + $chart1OffYDeltaPt

$co1.Width = $co1.Width + $chart1WidthDelta
$co1.Height = $co1.Height + $chart1HeightDelta

# Chart 2 ("from" corner shifts from col27 to col31, "to" corner moves
# from col43/row36 to col62/row38, 1-based).
$chart2LeftGain = ColsWidthPt 27 30
$chart2LeftDeltaPt = $chart2LeftGain + 0

$chart2ToColGain = ColsWidthPt 43 61
$chart2OffXDelta = 419100 - 281516
$chart2OffXDeltaPt = $chart2OffXDelta / $EMU_PER_PT
$chart2ToXDeltaPt = $chart2ToColGain + $chart2OffXDeltaPt

$chart2ToRowGain = RowsHeightPt 36 37
$chart2OffYDelta = 101600 - 152400
$chart2OffYDeltaPt = $chart2OffYDelta / $EMU_PER_PT
$chart2HeightDelta = $chart2ToRowGain + $chart2OffYDeltaPt

$co2.Left = $co2.Left + $chart2LeftDeltaPt
$co2.Width = $co2.Width + ($chart2ToXDeltaPt - $chart2LeftDeltaPt)
$co2.Height = $co2.Height + $chart2HeightDelta

# ---------------------------------------------------------------------
# 7. Restore the active selection.
# ---------------------------------------------------------------------
$ws.Range("D40").Select()
